$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# There is a gap in the data: row 16 is empty, with "ChatVRM" sitting at
# A17. Add the new "Flowise" dataset into that empty row 16, keeping
# "ChatVRM" where it is at row 17.
$ws.Range("A16").Value = "Flowise"

# Match the author's resulting selection (A16 instead of the old A17).
$ws.Range("A16").Select()
